# Update summary report in Excel format with latest data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Force text storage so numeric-looking strings (prices, SKUs) are not
    # auto-converted to numbers, then restore the default "Normal" style so
    # no stray number-format style gets attached to the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 3: HP OmniBook -> Dell Inspiron 15" (Core i7 / 16GB / 1TB)
Set-TextCell $ws.Range("A3") "Dell-  Inspiron 15"" Touch Screen Laptop  -  Intel Core i7 with 16GB Memory  -  1TB SSD  -  Black"
Set-TextCell $ws.Range("B3") "https://www.bestbuy.com/site/dell-inspiron-15-touch-screen-laptop-intel-core-i7-with-16gb-memory-1tb-ssd-black/6610571.p?skuId=6610571"
Set-TextCell $ws.Range("C3") "$649.99"
Set-TextCell $ws.Range("D3") "Rating 4.6 out of 5 stars with 81 reviews"
$ws.Range("E3").Value = 81
Set-TextCell $ws.Range("F3") "6610571"
Set-TextCell $ws.Range("G3") "i3530-7728BLK-PUS"

# Row 4: Lenovo Yoga -> Dell Inspiron 15" (Core i5 / 8GB / 512GB)
Set-TextCell $ws.Range("A4") "Dell-  Inspiron 15"" Touch Screen Laptop  -  Intel Core i5 with 8GB Memory  -  512GB SSD  -  Black"
Set-TextCell $ws.Range("B4") "https://www.bestbuy.com/site/dell-inspiron-15-touch-screen-laptop-intel-core-i5-with-8gb-memory-512gb-ssd-black/6610570.p?skuId=6610570"
Set-TextCell $ws.Range("C4") "$629.99"
Set-TextCell $ws.Range("D4") "Rating 4.7 out of 5 stars with 118 reviews"
$ws.Range("E4").Value = 118
Set-TextCell $ws.Range("F4") "6610570"
Set-TextCell $ws.Range("G4") "i3530-5623BLK-PUS"

# Row 5: HP OmniBook (Core 5) -> N/A (no longer available)
Set-TextCell $ws.Range("A5") "N/A"
Set-TextCell $ws.Range("B5") "N/A"
Set-TextCell $ws.Range("C5") "N/A"
Set-TextCell $ws.Range("D5") "N/A"
$ws.Range("E5").Value = 0
Set-TextCell $ws.Range("F5") "N/A"
Set-TextCell $ws.Range("G5") "N/A"

# Row 24: new trailing N/A row, extends used range to A1:G24
Set-TextCell $ws.Range("A24") "N/A"
Set-TextCell $ws.Range("B24") "N/A"
Set-TextCell $ws.Range("C24") "N/A"
Set-TextCell $ws.Range("D24") "N/A"
$ws.Range("E24").Value = 0
Set-TextCell $ws.Range("F24") "N/A"
Set-TextCell $ws.Range("G24") "N/A"
